$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case block (rows 16-17), mirroring the existing
# testSuccessfulIssuesFileDownLoad block in rows 11-12, for the new
# testSuccessfulIssuesFileDownLoadFiref test.

# Row 16: header row (bold/shaded/bordered, like row 11)
$ws.Range("A11:E11").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

# Row 17: data row (bordered, like row 12)
$ws.Range("A12:E12").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)

# Rows 18-20: trailing blank rows (columns B:E only), like rows 13-15
$ws.Range("B13:E13").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$ws.Range("B19:E19").PasteSpecial(-4122)
$ws.Range("B20:E20").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Cell values

$ws.Range("A16").Value = "testSuccessfulIssuesFileDownLoadFiref"
$ws.Range("B16").Value = "username"
$ws.Range("C16").Value = "password"
$ws.Range("D16").Value = "downloadPath"
$ws.Range("E16").Value = "filename"

$ws.Range("A17").Value = "testSuccessfulIssuesFileDownLoadFiref"
$ws.Range("B17").Value = "Richmond"
$ws.Range("C17").Value = "123456"
$ws.Range("D17").Value = "C:\\testfolder\\SeleniumDownloads"
$ws.Range("E17").Value = "Richmond"

$ws.Range("E17").Select()
